$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update header text (volume number, week-covering dates) ---
$ws.Range("A8").Value = "Volume 33   Number  5"
$ws.Range("C9").Value = "Report Covering the Week  1/26/2026  Through  2/1/2026"

# --- Simple same-style numeric value updates ---
$ws.Range("G14").Value = 2
$ws.Range("J14").Value = 2
$ws.Range("F15").Value = 2
$ws.Range("H15").Value = 100
$ws.Range("C16").Value = 1
$ws.Range("F16").Value = 8
$ws.Range("G16").Value = 6
$ws.Range("H16").Value = 33.333333333333
$ws.Range("I16").Value = 8
$ws.Range("J16").Value = 8
$ws.Range("K16").Value = 0
$ws.Range("L16").Value = -11.111111111111
$ws.Range("M16").Value = -80
$ws.Range("N16").Value = -91.752577319587
$ws.Range("C17").Value = 3
$ws.Range("E17").Value = -57.142857142857
$ws.Range("G17").Value = 21
$ws.Range("H17").Value = -42.857142857142
$ws.Range("I17").Value = 12
$ws.Range("J17").Value = 22
$ws.Range("K17").Value = -45.454545454545
$ws.Range("L17").Value = -55.555555555555
$ws.Range("M17").Value = 20
$ws.Range("N17").Value = -60
$ws.Range("C18").Value = 1
$ws.Range("D18").Value = 6
$ws.Range("E18").Value = -83.333333333333
$ws.Range("F18").Value = 11
$ws.Range("G18").Value = 26
$ws.Range("H18").Value = -57.692307692307
$ws.Range("I18").Value = 13
$ws.Range("J18").Value = 29
$ws.Range("K18").Value = -55.172413793103
$ws.Range("L18").Value = -48
$ws.Range("M18").Value = -55.172413793103
$ws.Range("N18").Value = -93.121693121693
$ws.Range("D19").Value = 14
$ws.Range("E19").Value = -42.857142857142
$ws.Range("F19").Value = 38
$ws.Range("G19").Value = 38
$ws.Range("H19").Value = 0
$ws.Range("I19").Value = 43
$ws.Range("J19").Value = 39
$ws.Range("K19").Value = 10.256410256410
$ws.Range("L19").Value = -4.444444444444
$ws.Range("M19").Value = -4.444444444444
$ws.Range("N19").Value = -6.521739130434
$ws.Range("C20").Value = 2
$ws.Range("D20").Value = 9
$ws.Range("E20").Value = -77.777777777777
$ws.Range("F20").Value = 22
$ws.Range("G20").Value = 28
$ws.Range("H20").Value = -21.428571428571
$ws.Range("I20").Value = 25
$ws.Range("J20").Value = 32
$ws.Range("K20").Value = -21.875
$ws.Range("L20").Value = -19.354838709677
$ws.Range("M20").Value = 25
$ws.Range("N20").Value = -94.252873563218
$ws.Range("C21").Value = 15
$ws.Range("D21").Value = 40
$ws.Range("E21").Value = -62.5
$ws.Range("F21").Value = 93
$ws.Range("G21").Value = 122
$ws.Range("H21").Value = -23.770491803278
$ws.Range("I21").Value = 103
$ws.Range("J21").Value = 133
$ws.Range("K21").Value = -22.556390977443
$ws.Range("L21").Value = -25.362318840579
$ws.Range("M21").Value = -29.452054794520
$ws.Range("N21").Value = -87.125
$ws.Range("E22").Value = 0
$ws.Range("F22").Value = 4
$ws.Range("G22").Value = 2
$ws.Range("H22").Value = 100
$ws.Range("I22").Value = 4
$ws.Range("J22").Value = 2
$ws.Range("K22").Value = 100
$ws.Range("L22").Value = 33.333333333333
$ws.Range("M22").Value = 100
$ws.Range("G23").Value = 2
$ws.Range("H23").Value = 100
$ws.Range("I23").Value = 4
$ws.Range("K23").Value = 100
$ws.Range("L23").Value = -33.333333333333
$ws.Range("M23").Value = 33.333333333333
$ws.Range("C24").Value = 17
$ws.Range("D24").Value = 30
$ws.Range("E24").Value = -43.333333333333
$ws.Range("F24").Value = 78
$ws.Range("G24").Value = 71
$ws.Range("H24").Value = 9.859154929577
$ws.Range("I24").Value = 96
$ws.Range("J24").Value = 75
$ws.Range("K24").Value = 28
$ws.Range("L24").Value = 4.347826086956
$ws.Range("M24").Value = 31.506849315068
$ws.Range("D25").Value = 17
$ws.Range("E25").Value = -100
$ws.Range("F25").Value = 14
$ws.Range("G25").Value = 26
$ws.Range("H25").Value = -46.153846153846
$ws.Range("I25").Value = 15
$ws.Range("J25").Value = 28
$ws.Range("K25").Value = -46.428571428571
$ws.Range("L25").Value = -64.285714285714
$ws.Range("C26").Value = 8
$ws.Range("E26").Value = -27.272727272727
$ws.Range("F26").Value = 28
$ws.Range("G26").Value = 45
$ws.Range("H26").Value = -37.777777777777
$ws.Range("I26").Value = 33
$ws.Range("J26").Value = 51
$ws.Range("K26").Value = -35.294117647058
$ws.Range("L26").Value = 13.793103448275
$ws.Range("M26").Value = -8.333333333333
$ws.Range("F27").Value = 2
$ws.Range("H27").Value = 100
$ws.Range("C28").Value = 1
$ws.Range("F28").Value = 6
$ws.Range("H28").Value = 100
$ws.Range("I28").Value = 7
$ws.Range("K28").Value = 133.333333333333
$ws.Range("L28").Value = 0
$ws.Range("G31").Value = 3
$ws.Range("J31").Value = 3

# --- Cells converting from text placeholder to numeric (reuse number format from a same-column/style donor) ---
$ws.Range("D14").Value = 1
$ws.Range("D14").NumberFormat = $ws.Range("G15").NumberFormat
$ws.Range("E14").Value = -100
$ws.Range("E14").NumberFormat = $ws.Range("H14").NumberFormat
$ws.Range("D16").Value = 3
$ws.Range("D16").NumberFormat = $ws.Range("G15").NumberFormat
$ws.Range("E16").Value = -66.666666666666
$ws.Range("E16").NumberFormat = $ws.Range("H14").NumberFormat
$ws.Range("C22").Value = 1
$ws.Range("C22").NumberFormat = $ws.Range("G15").NumberFormat
$ws.Range("C23").Value = 1
$ws.Range("C23").NumberFormat = $ws.Range("G15").NumberFormat
$ws.Range("D29").Value = 1
$ws.Range("D29").NumberFormat = $ws.Range("G15").NumberFormat
$ws.Range("E29").Value = -100
$ws.Range("E29").NumberFormat = $ws.Range("H14").NumberFormat
$ws.Range("G29").Value = 1
$ws.Range("G29").NumberFormat = $ws.Range("G15").NumberFormat
$ws.Range("H29").Value = -100
$ws.Range("H29").NumberFormat = $ws.Range("H14").NumberFormat
$ws.Range("J29").Value = 1
$ws.Range("J29").NumberFormat = $ws.Range("G15").NumberFormat
$ws.Range("K29").Value = -100
$ws.Range("K29").NumberFormat = $ws.Range("H14").NumberFormat
$ws.Range("D30").Value = 1
$ws.Range("D30").NumberFormat = $ws.Range("G15").NumberFormat
$ws.Range("E30").Value = -100
$ws.Range("E30").NumberFormat = $ws.Range("H14").NumberFormat
$ws.Range("G30").Value = 1
$ws.Range("G30").NumberFormat = $ws.Range("G15").NumberFormat
$ws.Range("H30").Value = -100
$ws.Range("H30").NumberFormat = $ws.Range("H14").NumberFormat
$ws.Range("J30").Value = 1
$ws.Range("J30").NumberFormat = $ws.Range("G15").NumberFormat
$ws.Range("K30").Value = -100
$ws.Range("K30").NumberFormat = $ws.Range("H14").NumberFormat

# --- Cells converting from numeric back to text placeholder (copy from a stable donor cell) ---
$ws.Range("C14").Copy($ws.Range("C25"))

Write-Host "Applied all cell updates"
